$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

$ws.Range('D2').Value = '64.483.40'
$ws.Range('E2').Value = '  -2.39%  '
$ws.Range('D3').Value = '3.180.77'
$ws.Range('E3').Value = '  -4.08%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextCell $ws.Range('D5') '572.43'
$ws.Range('E5').Value = '  -2.33%  '
Set-TextCell $ws.Range('D6') '169.53'
$ws.Range('E6').Value = '  -6.66%  '
Set-TextCell $ws.Range('D7') '0.610'
$ws.Range('E7').Value = '  -6.54%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '3.189.13'
$ws.Range('E9').Value = '  -3.85%  '
Set-TextCell $ws.Range('D10') '0.121'
$ws.Range('E10').Value = '  -3.82%  '
$ws.Range('E11').Value = '  +0.40%  '
Set-TextCell $ws.Range('D12') '0.390'
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('D13').Value = '3.733.64'
$ws.Range('E13').Value = '  -4.11%  '
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').Value = '64.517.87'
$ws.Range('E15').Value = '  -2.43%  '
Set-TextCell $ws.Range('D16') '25.35'
$ws.Range('E16').Value = '  -3.04%  '
$ws.Range('E17').Value = '  -3.55%  '
$ws.Range('D18').Value = '3.187.31'
$ws.Range('E18').Value = '  -3.73%  '
Set-TextCell $ws.Range('D19') '421.02'
$ws.Range('E19').Value = '  -0.94%  '
Set-TextCell $ws.Range('D20') '13.00'
$ws.Range('E20').Value = '  -0.87%  '
Set-TextCell $ws.Range('D21') '5.36'
$ws.Range('E21').Value = '  -3.17%  '
Set-TextCell $ws.Range('D22') '7.16'
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('E23').Value = '  -0.06%  '
Set-TextCell $ws.Range('D24') '70.32'
$ws.Range('E24').Value = '  -1.81%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  +2.39%  '
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('E28').Value = '  -7.60%  '
Set-TextCell $ws.Range('D29') '8.75'
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('E31').Value = '  -4.33%  '
Set-TextCell $ws.Range('D32') '21.80'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('E34').Value = '  -2.25%  '
Set-TextCell $ws.Range('D35') '6.37'
$ws.Range('E35').Value = '  -2.77%  '
Set-TextCell $ws.Range('D36') '157.09'
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('E37').Value = '  -4.07%  '
$ws.Range('E38').Value = '  -4.85%  '
$ws.Range('D39').Value = '2.710.31'
$ws.Range('E39').Value = '  -5.48%  '
Set-TextCell $ws.Range('D40') '1.71'
$ws.Range('E40').Value = '  -4.86%  '
Set-TextCell $ws.Range('D41') '24.34'
$ws.Range('E41').Value = '  -7.64%  '
$ws.Range('E42').Value = '  -1.44%  '
Set-TextCell $ws.Range('D43') '39.11'
$ws.Range('E43').Value = '  -1.79%  '
Set-TextCell $ws.Range('D44') '0.717'
$ws.Range('E44').Value = '  -5.44%  '
Set-TextCell $ws.Range('D45') '0.0622'
$ws.Range('E45').Value = '  -5.67%  '
Set-TextCell $ws.Range('D46') '5.51'
$ws.Range('E46').Value = '  -6.78%  '
$ws.Range('E47').Value = '  -2.89%  '
Set-TextCell $ws.Range('D48') '291.60'
$ws.Range('E48').Value = '  -6.69%  '
$ws.Range('E49').Value = '  -7.32%  '
Set-TextCell $ws.Range('D50') '0.0994'
$ws.Range('E50').Value = '  -5.76%  '
Set-TextCell $ws.Range('D51') '0.997'
$ws.Range('E51').Value = '  -0.26%  '
